$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# v33 fixed entropyTolerance and set to 0
$ws.Range("B5").Value = 0

# Update selection to B6 (matches post-edit cursor position in the diff)
$ws.Range("B6").Select()
